$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.328.46"
$ws.Range("E2").Value = "  +0.03%  "
$ws.Range("D3").Value = "1.876.93"
$ws.Range("E3").Value = "  +0.22%  "
$ws.Range("E4").Value = "  +0.02%  "
$cell = $ws.Range("D5")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.7117"
$cell.Style = $origStyle
$ws.Range("E5").Value = "  -0.03%  "
$cell = $ws.Range("D6")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "242.91"
$cell.Style = $origStyle
$ws.Range("E6").Value = "  +0.48%  "
$ws.Range("E7").Value = "  +0.00%  "
$cell = $ws.Range("D8")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.08013"
$cell.Style = $origStyle
$ws.Range("E8").Value = "  +3.19%  "
$cell = $ws.Range("D9")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.3155"
$cell.Style = $origStyle
$ws.Range("E9").Value = "  +1.60%  "
$ws.Range("E10").Value = "  -0.19%  "
$cell = $ws.Range("D11")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.08245"
$cell.Style = $origStyle
$ws.Range("E11").Value = "  -1.81%  "
$ws.Range("D12").Value = "1.896.99"
$ws.Range("E12").Value = "  +1.06%  "
$cell = $ws.Range("D13")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "5.253"
$cell.Style = $origStyle
$ws.Range("E13").Value = "  +0.34%  "
$cell = $ws.Range("D14")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "94.83"
$cell.Style = $origStyle
$ws.Range("E14").Value = "  +4.13%  "
$cell = $ws.Range("D15")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.7126"
$cell.Style = $origStyle
$ws.Range("E15").Value = "  +0.26%  "
$cell = $ws.Range("D16")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "6.384"
$cell.Style = $origStyle
$ws.Range("E16").Value = "  +5.37%  "
$cell = $ws.Range("D17")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.000008571"
$cell.Style = $origStyle
$ws.Range("E17").Value = "  +4.67%  "
$ws.Range("D18").Value = "29.353.34"
$ws.Range("E18").Value = "  +0.07%  "
$cell = $ws.Range("D19")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "243.94"
$cell.Style = $origStyle
$ws.Range("E19").Value = "  +1.87%  "
$ws.Range("D20").Value = "2.151.36"
$ws.Range("E20").Value = "  +1.39%  "
$cell = $ws.Range("D21")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "13.27"
$cell.Style = $origStyle
$ws.Range("E21").Value = "  +0.56%  "
$ws.Range("E22").Value = "  +0.05%  "
$cell = $ws.Range("D23")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "7.787"
$cell.Style = $origStyle
$ws.Range("E23").Value = "  +0.46%  "
$cell = $ws.Range("D24")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.001"
$cell.Style = $origStyle
$ws.Range("E24").Value = "  -0.02%  "
$cell = $ws.Range("D25")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.1564"
$cell.Style = $origStyle
$ws.Range("E25").Value = "  -2.05%  "
$cell = $ws.Range("D26")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "9.050"
$cell.Style = $origStyle
$ws.Range("E26").Value = "  +0.30%  "
$cell = $ws.Range("D27")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "162.61"
$cell.Style = $origStyle
$ws.Range("E27").Value = "  -0.05%  "
$cell = $ws.Range("D28")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "18.55"
$cell.Style = $origStyle
$ws.Range("E28").Value = "  +0.37%  "
$cell = $ws.Range("D29")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.501"
$cell.Style = $origStyle
$ws.Range("E29").Value = "  -0.45%  "
$cell = $ws.Range("D30")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "4.425"
$cell.Style = $origStyle
$ws.Range("E30").Value = "  +0.51%  "
$cell = $ws.Range("D31")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "4.311"
$cell.Style = $origStyle
$ws.Range("E31").Value = "  -0.18%  "
$cell = $ws.Range("D32")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.187"
$cell.Style = $origStyle
$ws.Range("E32").Value = "  -7.67%  "
$ws.Range("E33").Value = "  +1.68%  "
$cell = $ws.Range("D34")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.942"
$cell.Style = $origStyle
$ws.Range("E34").Value = "  +0.29%  "
$cell = $ws.Range("D35")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.7656"
$cell.Style = $origStyle
$ws.Range("E35").Value = "  +2.75%  "
$cell = $ws.Range("D36")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.181"
$cell.Style = $origStyle
$ws.Range("E36").Value = "  +0.52%  "
$ws.Range("E37").Value = "  -0.53%  "
$ws.Range("E38").Value = "  +0.30%  "
$ws.Range("D39").Value = "1.253.47"
$ws.Range("E39").Value = "  +2.85%  "
$cell = $ws.Range("D40")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.755"
$cell.Style = $origStyle
$ws.Range("E40").Value = "  +1.17%  "
$cell = $ws.Range("D41")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "6.487"
$cell.Style = $origStyle
$ws.Range("E41").Value = "  -0.38%  "
$cell = $ws.Range("D42")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.9177"
$cell.Style = $origStyle
$ws.Range("E42").Value = "  +3.54%  "
$cell = $ws.Range("D43")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "112.82"
$cell.Style = $origStyle
$ws.Range("E43").Value = "  +2.76%  "
$ws.Range("E44").Value = "  +2.44%  "
$ws.Range("E45").Value = "  +9.28%  "
$ws.Range("E46").Value = "  +0.02%  "
$ws.Range("D47").Value = "2.041.99"
$ws.Range("E47").Value = "  +1.13%  "
$cell = $ws.Range("D48")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.5222"
$cell.Style = $origStyle
$ws.Range("E48").Value = "  +0.53%  "
$cell = $ws.Range("D49")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "1.804"
$cell.Style = $origStyle
$ws.Range("E49").Value = "  +0.30%  "
$cell = $ws.Range("D50")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "9.474"
$cell.Style = $origStyle
$ws.Range("E50").Value = "  +1.34%  "
$cell = $ws.Range("D51")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.4359"
$cell.Style = $origStyle
$ws.Range("E51").Value = "  +1.07%  "
